$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.219.23'
$ws.Range("E2").Value = '  +2.80%  '

# Row 3
$ws.Range("D3").Value = '3.417.13'
$ws.Range("E3").Value = '  +2.18%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.20%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.40%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("E8").Value = '  +0.25%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.56'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.58%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.128'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.92%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.394'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.72%  '

# Row 12
$ws.Range("D12").Value = '3.998.08'
$ws.Range("E12").Value = '  +2.12%  '

# Row 13
$ws.Range("E13").Value = '  +2.05%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000180'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.62%  '

# Row 15
$ws.Range("D15").Value = '3.412.55'
$ws.Range("E15").Value = '  +2.14%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.51%  '

# Row 17
$ws.Range("D17").Value = '62.213.29'
$ws.Range("E17").Value = '  +2.64%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.24'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.92%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.21%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.40%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '396.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.79%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.571'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.39%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000132'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +14.90%  '

# Row 24
$ws.Range("D24").Value = '3.551.02'
$ws.Range("E24").Value = '  +2.22%  '

# Row 25
$ws.Range("E25").Value = '  +0.09%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.51%  '

# Row 27
$ws.Range("B27").Value = 'Fetch.AI'
$ws.Range("C27").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.79%  '

# Row 28
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.98%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.00%  '

# Row 30
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.163'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.10%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.36%  '

# Row 32
$ws.Range("E32").Value = '  +2.39%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.64'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.34%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.58%  '

# Row 36
$ws.Range("D36").Value = '3.447.91'
$ws.Range("E36").Value = '  +2.29%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.60'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.07%  '

# Row 38
$ws.Range("E38").Value = '  +0.78%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.30%  '

# Row 40
$ws.Range("E40").Value = '  +1.84%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +13.66%  '

# Row 42
$ws.Range("E42").Value = '  +4.61%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.788'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.12%  '

# Row 44
$ws.Range("E44").Value = '  -0.03%  '

# Row 45
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.65%  '

# Row 46
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.37%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.42%  '

# Row 48
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.57'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.85%  '

# Row 49
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.96'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.26%  '

# Row 50
$ws.Range("D50").Value = '2.371.40'
$ws.Range("E50").Value = '  +9.36%  '

# Row 51
$ws.Range("E51").Value = '  -1.24%  '
